$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric attribute columns (D = dano, E = cadencia_tiro) ---
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = 200

$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 250

$ws.Range("D4").Value = 120
$ws.Range("E4").Value = 100

$ws.Range("D5").Value = 200
$ws.Range("E5").Value = 100

$ws.Range("D6").Value = 1000
$ws.Range("E6").Value = 500

$ws.Range("D7").Value = 500
$ws.Range("E7").Value = 150

$ws.Range("D8").Value = 250
$ws.Range("E8").Value = 20

$ws.Range("D9").Value = 2000
$ws.Range("E9").Value = 600

$ws.Range("D10").Value = 2000
$ws.Range("E10").Value = 250

# --- Update column G (tiro_automatico) so every row becomes the text "false" ---
# A direct string assignment of "false"/"true" gets auto-coerced by Excel into a
# Boolean cell (t="b"), which does not match the original text-string ("true"/"false")
# cells already present in the workbook. To keep the cell as a shared-string text
# value we build a text "false" via a TEXT() formula, copy it, and paste-special
# only the resulting value into each target cell (this preserves the text type).
$ws.Range("Z1").Formula = '=TEXT(0,"\f\a\l\s\e")'

$ws.Range("Z1").Copy()
$ws.Range("G2").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G3").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G4").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G5").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G6").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G7").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G8").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G9").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G10").PasteSpecial(-4163)

# Clean up the helper cell used to stage the text "false" value.
$ws.Range("Z1").Clear()

# --- Update the active selection shown in the sheet view ---
$ws.Range("I2").Select()
